{"js": "// Update the division-fact answers in the single table of the document.\n// Only the non-empty data rows (0, 4, 8, 12, 16 \u2014 0-indexed) contain text;\n// the intervening rows are blank spacer rows and are left untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"91\u00f79=10, 1\", \"68\u00f76=11, 2\", \"81\u00f75=16, 1\", \"82\u00f74=20, 2\", \"51\u00f79=5, 6\"],\n  [\"53\u00f75=10, 3\", \"10\u00f75=2, 0\", \"92\u00f76=15, 2\", \"15\u00f72=7, 1\", \"92\u00f72=46, 0\"],\n  [\"96\u00f75=19, 1\", \"22\u00f74=5, 2\", \"37\u00f74=9, 1\", \"12\u00f72=6, 0\", \"69\u00f73=23, 0\"],\n  [\"17\u00f78=2, 1\", \"90\u00f76=15, 0\", \"33\u00f77=4, 5\", \"56\u00f75=11, 1\", \"23\u00f72=11, 1\"],\n  [\"22\u00f72=11, 0\", \"47\u00f72=23, 1\", \"97\u00f72=48, 1\", \"26\u00f74=6, 2\", \"46\u00f75=9, 1\"],\n];\n\nconst dataRows = [0, 4, 8, 12, 16];\n\nfor (let i = 0; i < dataRows.length; i++) {\n  const rowIndex = dataRows[i];\n  const rowValues = newValues[i];\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-fact answers in the single table of the document.\n# Only the non-empty data rows (1, 5, 9, 13, 17 \u2014 1-indexed) contain text;\n# the intervening rows are blank spacer rows and are left untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"91\u00f79=10, 1\", \"68\u00f76=11, 2\", \"81\u00f75=16, 1\", \"82\u00f74=20, 2\", \"51\u00f79=5, 6\"),\n    @(\"53\u00f75=10, 3\", \"10\u00f75=2, 0\",  \"92\u00f76=15, 2\", \"15\u00f72=7, 1\",  \"92\u00f72=46, 0\"),\n    @(\"96\u00f75=19, 1\", \"22\u00f74=5, 2\",  \"37\u00f74=9, 1\",  \"12\u00f72=6, 0\",  \"69\u00f73=23, 0\"),\n    @(\"17\u00f78=2, 1\",  \"90\u00f76=15, 0\", \"33\u00f77=4, 5\",  \"56\u00f75=11, 1\", \"23\u00f72=11, 1\"),\n    @(\"22\u00f72=11, 0\", \"47\u00f72=23, 1\", \"97\u00f72=48, 1\", \"26\u00f74=6, 2\",  \"46\u00f75=9, 1\")\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $dataRows.Count; $i++) {\n    $rowIndex = $dataRows[$i]\n    $rowValues = $newValues[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
